$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 corresponds to the D1 LED part. Restore its color back to Red,
# along with the matching manufacturer/part number for the red LED variant.
$ws.Range("B16").Value = "Red"
$ws.Range("E16").Value = "ROHM"
$ws.Range("F16").Value = "SML-D12U8WT86C"

# Match the author's final selection/active cell.
$ws.Range("F16").Select()
